$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44172
$ws.Range("M3").Value = 90
$ws.Range("N3").Value = 8500
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 8806
$ws.Range("S3").Value = 629

# Row 4
$ws.Range("D4").Value = 44229
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 11000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 11364
$ws.Range("S4").Value = 812

# Row 5
$ws.Range("D5").Value = 44181
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 9000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 9462
$ws.Range("S5").Value = 676

# Row 6
$ws.Range("D6").Value = 44216
$ws.Range("M6").Value = 55
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11545
$ws.Range("S6").Value = 825

# Row 8
$ws.Range("D8").Value = 44253
$ws.Range("M8").Value = 90
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 12667
$ws.Range("S8").Value = 905
